$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2021 -Month 11 -Day 18 -Hour 0 -Minute 0 -Second 0

# Insert 3 new rows before row 180 (weekly data update), shifting existing rows down.
$ws.Rows("180:182").Insert()

# Row 180
$ws.Cells.Item(180, 1).Value = 6
$ws.Cells.Item(180, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(180, 3).Value = "Metropolitana"
$ws.Cells.Item(180, 4).Value = $newDate
$ws.Cells.Item(180, 5).Value = 13
$ws.Cells.Item(180, 6).Value = 100112052
$ws.Cells.Item(180, 7).Value = "Albahaca"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 400
$ws.Cells.Item(180, 11).Value = 5000
$ws.Cells.Item(180, 12).Value = 5500
$ws.Cells.Item(180, 13).Value = 5188
$ws.Cells.Item(180, 14).Value = "`$/docena de matas"
$ws.Cells.Item(180, 15).Value = "Región Metropolitana"
$ws.Cells.Item(180, 16).Value = 865
$ws.Cells.Item(180, 17).Value = 6
$ws.Cells.Item(180, 18).Value = "Hortaliza"

# Row 181
$ws.Cells.Item(181, 1).Value = 6
$ws.Cells.Item(181, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(181, 3).Value = "Metropolitana"
$ws.Cells.Item(181, 4).Value = $newDate
$ws.Cells.Item(181, 5).Value = 13
$ws.Cells.Item(181, 6).Value = 100112052
$ws.Cells.Item(181, 7).Value = "Albahaca"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 160
$ws.Cells.Item(181, 11).Value = 4500
$ws.Cells.Item(181, 12).Value = 5000
$ws.Cells.Item(181, 13).Value = 4781
$ws.Cells.Item(181, 14).Value = "`$/paquete"
$ws.Cells.Item(181, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(181, 16).Value = 4781
$ws.Cells.Item(181, 17).Value = 1
$ws.Cells.Item(181, 18).Value = "Hortaliza"

# Row 182
$ws.Cells.Item(182, 1).Value = 6
$ws.Cells.Item(182, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(182, 3).Value = "Metropolitana"
$ws.Cells.Item(182, 4).Value = $newDate
$ws.Cells.Item(182, 5).Value = 13
$ws.Cells.Item(182, 6).Value = 100112052
$ws.Cells.Item(182, 7).Value = "Albahaca"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Segunda"
$ws.Cells.Item(182, 10).Value = 100
$ws.Cells.Item(182, 11).Value = 4500
$ws.Cells.Item(182, 12).Value = 4500
$ws.Cells.Item(182, 13).Value = 4500
$ws.Cells.Item(182, 14).Value = "`$/docena de matas"
$ws.Cells.Item(182, 15).Value = "Región Metropolitana"
$ws.Cells.Item(182, 16).Value = 750
$ws.Cells.Item(182, 17).Value = 6
$ws.Cells.Item(182, 18).Value = "Hortaliza"
